{"js": "// Change the highlight color of the \"reminder\" user-story bullet from\n// yellow to green (paragraph mark + every run in the paragraph).\nconst searchText = \"set certain expenses to appear as a reminder on my home page\";\nconst results = context.document.body.search(searchText, { matchCase: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target paragraph (reminder user story).\");\n}\n\nconst targetParagraph = results.items[0].paragraphs.getFirst();\n\n// Setting highlightColor on the paragraph's font updates both the\n// paragraph mark run properties (w:pPr/w:rPr) and every run inside it,\n// which is exactly what's needed to flip the whole bullet yellow -> green.\ntargetParagraph.font.highlightColor = \"green\";\n\nawait context.sync();\n", "ps1": "# Change the highlight color of the \"reminder\" user-story bullet from\n# yellow to green (paragraph mark + every run in the paragraph).\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*set certain expenses to appear as a reminder on my home page*\") {\n        # Font.HighlightColorIndex (vs. Range.HighlightColorIndex) applies to\n        # the paragraph mark run properties (w:pPr/w:rPr) as well as every\n        # run inside the paragraph - matching the whole bullet turning green.\n        $p.Range.Font.HighlightColorIndex = 4\n    }\n}\n"}
